# "Progress on Sample 5"
# Shift the measured-data table up by one row (the new Sample 5 reading
# now occupies row 1) and refine the conversion constant used in column B
# from 3.14 to 3.141596.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values (same samples, now starting at row 1 instead of row 2)
$ws.Range("A1").Value = 140
$ws.Range("A2").Value = 67
$ws.Range("A3").Value = 50.5
$ws.Range("A4").Value = 12.5
$ws.Range("A5").Value = 146.30000000000001
$ws.Range("A6").Value = 92.8

# Drop the now-unused trailing row (was row 7 before the shift)
$ws.Range("A7:B7").ClearContents()

# Recreate column B with the updated divisor constant, filling the shared
# formula block starting at row 2 first so the fill-anchor lands on B2
# (matching row1 getting its own, non-shared copy of the formula).
$ws.Range("B2:B6").Formula = "=A2/3.141596"
$ws.Range("B1").Formula = "=A1/3.141596"

# Update the active selection to follow the data (was B5, now B4)
$ws.Range("B4").Select()
